$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.966.19"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "1.640.84"
$ws.Range("E3").Value = "  +0.50%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "212.95"
$ws.Range("E5").Value = "  +0.42%  "
$ws.Range("E6").Value = "  +0.16%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "23.58"
$ws.Range("E8").Value = "  +1.75%  "
$ws.Range("E9").Value = "  -1.47%  "
$ws.Range("E10").Value = "  +0.41%  "
$ws.Range("E11").Value = "  +2.16%  "
$ws.Range("D12").Value = "1.873.14"
$ws.Range("E12").Value = "  +0.47%  "
$ws.Range("D13").Value = "1.640.62"
$ws.Range("E13").Value = "  +0.37%  "
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").Value = "0.575"
$ws.Range("E15").Value = "  +4.05%  "
$ws.Range("D16").Value = "65.92"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").Value = "27.962.65"
$ws.Range("E17").Value = "  +1.13%  "
$ws.Range("D18").Value = "233.96"
$ws.Range("E18").Value = "  +1.77%  "
$ws.Range("E19").Value = "  +0.68%  "
$ws.Range("D20").Value = "7.63"
$ws.Range("E20").Value = "  +0.73%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "10.71"
$ws.Range("E22").Value = "  +0.02%  "
$ws.Range("E23").Value = "  +0.70%  "
$ws.Range("E24").Value = "  -2.28%  "
$ws.Range("D25").Value = "151.03"
$ws.Range("E25").Value = "  +1.54%  "
$ws.Range("E26").Value = "  +1.51%  "
$ws.Range("D27").Value = "15.68"
$ws.Range("E27").Value = "  +0.66%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E31").Value = "  +0.10%  "
$ws.Range("E32").Value = "  +2.03%  "
$ws.Range("B33").Value = "Maker"
$ws.Range("C33").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D33").Value = "1.424.99"
$ws.Range("E33").Value = "  -3.17%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D34").Value = "3.12"
$ws.Range("E34").Value = "  +1.24%  "
$ws.Range("E35").Value = "  +2.26%  "
$ws.Range("D36").Value = "2.36"
$ws.Range("E36").Value = "  +1.21%  "
$ws.Range("E37").Value = "  +1.28%  "
$ws.Range("D38").Value = "0.882"
$ws.Range("E38").Value = "  +0.57%  "
$ws.Range("E39").Value = "  -0.18%  "
$ws.Range("E40").Value = "  -4.31%  "
$ws.Range("E41").Value = "  +0.82%  "
$ws.Range("E42").Value = "  -0.10%  "
$ws.Range("E43").Value = "  +6.80%  "
$ws.Range("E44").Value = "  -2.03%  "
$ws.Range("D45").Value = "5.51"
$ws.Range("E46").Value = "  +0.07%  "
$ws.Range("D47").Value = "1.782.07"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("E48").Value = "  +0.46%  "
$ws.Range("E49").Value = "  +1.22%  "
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").Value = "7.63"
$ws.Range("E51").Value = "  -1.03%  "
